$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: numeric corrections to existing rows (no district-name change) ---
$ws.Cells.Item(2539, 3).Value = 99
$ws.Cells.Item(2539, 5).Value = 102
$ws.Cells.Item(2605, 3).Value = 395
$ws.Cells.Item(2605, 5).Value = 401
$ws.Cells.Item(2909, 3).Value = 294
$ws.Cells.Item(2909, 5).Value = 308
$ws.Cells.Item(2911, 3).Value = 855
$ws.Cells.Item(2911, 5).Value = 873

# row 2939 has C, D, and E all corrected
$ws.Cells.Item(2939, 3).Value = 941
$ws.Cells.Item(2939, 4).Value = 91
$ws.Cells.Item(2939, 5).Value = 1032

# row 2971 (Okres Zarnovica @ 44173) values corrected
$ws.Cells.Item(2971, 3).Value = 1575
$ws.Cells.Item(2971, 4).Value = 132
$ws.Cells.Item(2971, 5).Value = 1707

# --- Part 2: rewrite rows 2974-3027 (date 44174, now reordered/extended) ---
# and append rows 3028-3100 (rest of 44174 + new date 44175 block) ---
$ws.Cells.Item(2974, 1).Value = 44174
$ws.Cells.Item(2974, 2).Value = "Okres Bánovce nad Bebravou"
$ws.Cells.Item(2974, 3).Value = 189
$ws.Cells.Item(2974, 4).Value = 27
$ws.Cells.Item(2974, 5).Value = 216
$ws.Cells.Item(2975, 1).Value = 44174
$ws.Cells.Item(2975, 2).Value = "Okres Banská Bystrica"
$ws.Cells.Item(2975, 3).Value = 606
$ws.Cells.Item(2975, 4).Value = 52
$ws.Cells.Item(2975, 5).Value = 658
$ws.Cells.Item(2976, 1).Value = 44174
$ws.Cells.Item(2976, 2).Value = "Okres Banská Štiavnica"
$ws.Cells.Item(2976, 3).Value = 15
$ws.Cells.Item(2976, 4).Value = 3
$ws.Cells.Item(2976, 5).Value = 18
$ws.Cells.Item(2977, 1).Value = 44174
$ws.Cells.Item(2977, 2).Value = "Okres Bardejov"
$ws.Cells.Item(2977, 3).Value = 144
$ws.Cells.Item(2977, 4).Value = 7
$ws.Cells.Item(2977, 5).Value = 151
$ws.Cells.Item(2978, 1).Value = 44174
$ws.Cells.Item(2978, 2).Value = "Okres Bratislava I"
$ws.Cells.Item(2978, 3).Value = 215
$ws.Cells.Item(2978, 4).Value = 7
$ws.Cells.Item(2978, 5).Value = 222
$ws.Cells.Item(2979, 1).Value = 44174
$ws.Cells.Item(2979, 2).Value = "Okres Bratislava II"
$ws.Cells.Item(2979, 3).Value = 1001
$ws.Cells.Item(2979, 4).Value = 22
$ws.Cells.Item(2979, 5).Value = 1023
$ws.Cells.Item(2980, 1).Value = 44174
$ws.Cells.Item(2980, 2).Value = "Okres Bratislava III"
$ws.Cells.Item(2980, 3).Value = 868
$ws.Cells.Item(2980, 4).Value = 30
$ws.Cells.Item(2980, 5).Value = 898
$ws.Cells.Item(2981, 1).Value = 44174
$ws.Cells.Item(2981, 2).Value = "Okres Bratislava V"
$ws.Cells.Item(2981, 3).Value = 695
$ws.Cells.Item(2981, 4).Value = 26
$ws.Cells.Item(2981, 5).Value = 721
$ws.Cells.Item(2982, 1).Value = 44174
$ws.Cells.Item(2982, 2).Value = "Okres Čadca"
$ws.Cells.Item(2982, 3).Value = 452
$ws.Cells.Item(2982, 4).Value = 21
$ws.Cells.Item(2982, 5).Value = 473
$ws.Cells.Item(2983, 1).Value = 44174
$ws.Cells.Item(2983, 2).Value = "Okres Detva"
$ws.Cells.Item(2983, 3).Value = 211
$ws.Cells.Item(2983, 4).Value = 18
$ws.Cells.Item(2983, 5).Value = 229
$ws.Cells.Item(2984, 1).Value = 44174
$ws.Cells.Item(2984, 2).Value = "Okres Dolný Kubín"
$ws.Cells.Item(2984, 3).Value = 147
$ws.Cells.Item(2984, 4).Value = 12
$ws.Cells.Item(2984, 5).Value = 159
$ws.Cells.Item(2985, 1).Value = 44174
$ws.Cells.Item(2985, 2).Value = "Okres Dunajská Streda"
$ws.Cells.Item(2985, 3).Value = 532
$ws.Cells.Item(2985, 4).Value = 34
$ws.Cells.Item(2985, 5).Value = 566
$ws.Cells.Item(2986, 1).Value = 44174
$ws.Cells.Item(2986, 2).Value = "Okres Galanta"
$ws.Cells.Item(2986, 3).Value = 550
$ws.Cells.Item(2986, 4).Value = 36
$ws.Cells.Item(2986, 5).Value = 586
$ws.Cells.Item(2987, 1).Value = 44174
$ws.Cells.Item(2987, 2).Value = "Okres Gelnica"
$ws.Cells.Item(2987, 3).Value = 168
$ws.Cells.Item(2987, 4).Value = 7
$ws.Cells.Item(2987, 5).Value = 175
$ws.Cells.Item(2988, 1).Value = 44174
$ws.Cells.Item(2988, 2).Value = "Okres Hlohovec"
$ws.Cells.Item(2988, 3).Value = 25
$ws.Cells.Item(2988, 4).Value = 0
$ws.Cells.Item(2988, 5).Value = 25
$ws.Cells.Item(2989, 1).Value = 44174
$ws.Cells.Item(2989, 2).Value = "Okres Humenné"
$ws.Cells.Item(2989, 3).Value = 403
$ws.Cells.Item(2989, 4).Value = 15
$ws.Cells.Item(2989, 5).Value = 418
$ws.Cells.Item(2990, 1).Value = 44174
$ws.Cells.Item(2990, 2).Value = "Okres Ilava"
$ws.Cells.Item(2990, 3).Value = 288
$ws.Cells.Item(2990, 4).Value = 36
$ws.Cells.Item(2990, 5).Value = 324
$ws.Cells.Item(2991, 1).Value = 44174
$ws.Cells.Item(2991, 2).Value = "Okres Komárno"
$ws.Cells.Item(2991, 3).Value = 735
$ws.Cells.Item(2991, 4).Value = 18
$ws.Cells.Item(2991, 5).Value = 753
$ws.Cells.Item(2992, 1).Value = 44174
$ws.Cells.Item(2992, 2).Value = "Okres Košice I"
$ws.Cells.Item(2992, 3).Value = 533
$ws.Cells.Item(2992, 4).Value = 32
$ws.Cells.Item(2992, 5).Value = 565
$ws.Cells.Item(2993, 1).Value = 44174
$ws.Cells.Item(2993, 2).Value = "Okres Košice II"
$ws.Cells.Item(2993, 3).Value = 580
$ws.Cells.Item(2993, 4).Value = 2
$ws.Cells.Item(2993, 5).Value = 582
$ws.Cells.Item(2994, 1).Value = 44174
$ws.Cells.Item(2994, 2).Value = "Okres Košice IV"
$ws.Cells.Item(2994, 3).Value = 247
$ws.Cells.Item(2994, 4).Value = 11
$ws.Cells.Item(2994, 5).Value = 258
$ws.Cells.Item(2995, 1).Value = 44174
$ws.Cells.Item(2995, 2).Value = "Okres Krupina"
$ws.Cells.Item(2995, 3).Value = 0
$ws.Cells.Item(2995, 4).Value = 0
$ws.Cells.Item(2995, 5).Value = 0
$ws.Cells.Item(2996, 1).Value = 44174
$ws.Cells.Item(2996, 2).Value = "Okres Levice"
$ws.Cells.Item(2996, 3).Value = 1491
$ws.Cells.Item(2996, 4).Value = 95
$ws.Cells.Item(2996, 5).Value = 1586
$ws.Cells.Item(2997, 1).Value = 44174
$ws.Cells.Item(2997, 2).Value = "Okres Levoča"
$ws.Cells.Item(2997, 3).Value = 229
$ws.Cells.Item(2997, 4).Value = 19
$ws.Cells.Item(2997, 5).Value = 248
$ws.Cells.Item(2998, 1).Value = 44174
$ws.Cells.Item(2998, 2).Value = "Okres Liptovský Mikuláš"
$ws.Cells.Item(2998, 3).Value = 599
$ws.Cells.Item(2998, 4).Value = 44
$ws.Cells.Item(2998, 5).Value = 643
$ws.Cells.Item(2999, 1).Value = 44174
$ws.Cells.Item(2999, 2).Value = "Okres Lučenec"
$ws.Cells.Item(2999, 3).Value = 456
$ws.Cells.Item(2999, 4).Value = 22
$ws.Cells.Item(2999, 5).Value = 478
$ws.Cells.Item(3000, 1).Value = 44174
$ws.Cells.Item(3000, 2).Value = "Okres Malacky"
$ws.Cells.Item(3000, 3).Value = 154
$ws.Cells.Item(3000, 4).Value = 12
$ws.Cells.Item(3000, 5).Value = 166
$ws.Cells.Item(3001, 1).Value = 44174
$ws.Cells.Item(3001, 2).Value = "Okres Martin"
$ws.Cells.Item(3001, 3).Value = 692
$ws.Cells.Item(3001, 4).Value = 23
$ws.Cells.Item(3001, 5).Value = 715
$ws.Cells.Item(3002, 1).Value = 44174
$ws.Cells.Item(3002, 2).Value = "Okres Medzilaborce"
$ws.Cells.Item(3002, 3).Value = 59
$ws.Cells.Item(3002, 4).Value = 1
$ws.Cells.Item(3002, 5).Value = 60
$ws.Cells.Item(3003, 1).Value = 44174
$ws.Cells.Item(3003, 2).Value = "Okres Michalovce"
$ws.Cells.Item(3003, 3).Value = 1118
$ws.Cells.Item(3003, 4).Value = 55
$ws.Cells.Item(3003, 5).Value = 1173
$ws.Cells.Item(3004, 1).Value = 44174
$ws.Cells.Item(3004, 2).Value = "Okres Myjava"
$ws.Cells.Item(3004, 3).Value = 347
$ws.Cells.Item(3004, 4).Value = 27
$ws.Cells.Item(3004, 5).Value = 374
$ws.Cells.Item(3005, 1).Value = 44174
$ws.Cells.Item(3005, 2).Value = "Okres Námestovo"
$ws.Cells.Item(3005, 3).Value = 126
$ws.Cells.Item(3005, 4).Value = 6
$ws.Cells.Item(3005, 5).Value = 132
$ws.Cells.Item(3006, 1).Value = 44174
$ws.Cells.Item(3006, 2).Value = "Okres Nitra"
$ws.Cells.Item(3006, 3).Value = 321
$ws.Cells.Item(3006, 4).Value = 22
$ws.Cells.Item(3006, 5).Value = 343
$ws.Cells.Item(3007, 1).Value = 44174
$ws.Cells.Item(3007, 2).Value = "Okres Nové Mesto nad Váhom"
$ws.Cells.Item(3007, 3).Value = 511
$ws.Cells.Item(3007, 4).Value = 47
$ws.Cells.Item(3007, 5).Value = 558
$ws.Cells.Item(3008, 1).Value = 44174
$ws.Cells.Item(3008, 2).Value = "Okres Nové Zámky"
$ws.Cells.Item(3008, 3).Value = 452
$ws.Cells.Item(3008, 4).Value = 201
$ws.Cells.Item(3008, 5).Value = 653
$ws.Cells.Item(3009, 1).Value = 44174
$ws.Cells.Item(3009, 2).Value = "Okres Partizánske"
$ws.Cells.Item(3009, 3).Value = 271
$ws.Cells.Item(3009, 4).Value = 33
$ws.Cells.Item(3009, 5).Value = 304
$ws.Cells.Item(3010, 1).Value = 44174
$ws.Cells.Item(3010, 2).Value = "Okres Pezinok"
$ws.Cells.Item(3010, 3).Value = 588
$ws.Cells.Item(3010, 4).Value = 45
$ws.Cells.Item(3010, 5).Value = 633
$ws.Cells.Item(3011, 1).Value = 44174
$ws.Cells.Item(3011, 2).Value = "Okres Piešťany"
$ws.Cells.Item(3011, 3).Value = 587
$ws.Cells.Item(3011, 4).Value = 38
$ws.Cells.Item(3011, 5).Value = 625
$ws.Cells.Item(3012, 1).Value = 44174
$ws.Cells.Item(3012, 2).Value = "Okres Poprad"
$ws.Cells.Item(3012, 3).Value = 455
$ws.Cells.Item(3012, 4).Value = 53
$ws.Cells.Item(3012, 5).Value = 508
$ws.Cells.Item(3013, 1).Value = 44174
$ws.Cells.Item(3013, 2).Value = "Okres Považská Bystrica"
$ws.Cells.Item(3013, 3).Value = 284
$ws.Cells.Item(3013, 4).Value = 34
$ws.Cells.Item(3013, 5).Value = 318
$ws.Cells.Item(3014, 1).Value = 44174
$ws.Cells.Item(3014, 2).Value = "Okres Prešov"
$ws.Cells.Item(3014, 3).Value = 1783
$ws.Cells.Item(3014, 4).Value = 95
$ws.Cells.Item(3014, 5).Value = 1878
$ws.Cells.Item(3015, 1).Value = 44174
$ws.Cells.Item(3015, 2).Value = "Okres Prievidza"
$ws.Cells.Item(3015, 3).Value = 545
$ws.Cells.Item(3015, 4).Value = 58
$ws.Cells.Item(3015, 5).Value = 603
$ws.Cells.Item(3016, 1).Value = 44174
$ws.Cells.Item(3016, 2).Value = "Okres Púchov"
$ws.Cells.Item(3016, 3).Value = 3
$ws.Cells.Item(3016, 4).Value = 0
$ws.Cells.Item(3016, 5).Value = 3
$ws.Cells.Item(3017, 1).Value = 44174
$ws.Cells.Item(3017, 2).Value = "Okres Revúca"
$ws.Cells.Item(3017, 3).Value = 326
$ws.Cells.Item(3017, 4).Value = 3
$ws.Cells.Item(3017, 5).Value = 329
$ws.Cells.Item(3018, 1).Value = 44174
$ws.Cells.Item(3018, 2).Value = "Okres Rimavská Sobota"
$ws.Cells.Item(3018, 3).Value = 285
$ws.Cells.Item(3018, 4).Value = 8
$ws.Cells.Item(3018, 5).Value = 293
$ws.Cells.Item(3019, 1).Value = 44174
$ws.Cells.Item(3019, 2).Value = "Okres Rožňava"
$ws.Cells.Item(3019, 3).Value = 204
$ws.Cells.Item(3019, 4).Value = 4
$ws.Cells.Item(3019, 5).Value = 208
$ws.Cells.Item(3020, 1).Value = 44174
$ws.Cells.Item(3020, 2).Value = "Okres Ružomberok"
$ws.Cells.Item(3020, 3).Value = 410
$ws.Cells.Item(3020, 4).Value = 36
$ws.Cells.Item(3020, 5).Value = 446
$ws.Cells.Item(3021, 1).Value = 44174
$ws.Cells.Item(3021, 2).Value = "Okres Sabinov"
$ws.Cells.Item(3021, 3).Value = 244
$ws.Cells.Item(3021, 4).Value = 16
$ws.Cells.Item(3021, 5).Value = 260
$ws.Cells.Item(3022, 1).Value = 44174
$ws.Cells.Item(3022, 2).Value = "Okres Senec"
$ws.Cells.Item(3022, 3).Value = 1081
$ws.Cells.Item(3022, 4).Value = 39
$ws.Cells.Item(3022, 5).Value = 1120
$ws.Cells.Item(3023, 1).Value = 44174
$ws.Cells.Item(3023, 2).Value = "Okres Senica"
$ws.Cells.Item(3023, 3).Value = 243
$ws.Cells.Item(3023, 4).Value = 30
$ws.Cells.Item(3023, 5).Value = 273
$ws.Cells.Item(3024, 1).Value = 44174
$ws.Cells.Item(3024, 2).Value = "Okres Skalica"
$ws.Cells.Item(3024, 3).Value = 569
$ws.Cells.Item(3024, 4).Value = 13
$ws.Cells.Item(3024, 5).Value = 582
$ws.Cells.Item(3025, 1).Value = 44174
$ws.Cells.Item(3025, 2).Value = "Okres Snina"
$ws.Cells.Item(3025, 3).Value = 13
$ws.Cells.Item(3025, 4).Value = 3
$ws.Cells.Item(3025, 5).Value = 16
$ws.Cells.Item(3026, 1).Value = 44174
$ws.Cells.Item(3026, 2).Value = "Okres Sobrance"
$ws.Cells.Item(3026, 3).Value = 285
$ws.Cells.Item(3026, 4).Value = 11
$ws.Cells.Item(3026, 5).Value = 296
$ws.Cells.Item(3027, 1).Value = 44174
$ws.Cells.Item(3027, 2).Value = "Okres Spišská Nová Ves"
$ws.Cells.Item(3027, 3).Value = 339
$ws.Cells.Item(3027, 4).Value = 31
$ws.Cells.Item(3027, 5).Value = 370
$ws.Cells.Item(3028, 1).Value = 44174
$ws.Cells.Item(3028, 2).Value = "Okres Stará Ľubovňa"
$ws.Cells.Item(3028, 3).Value = 208
$ws.Cells.Item(3028, 4).Value = 13
$ws.Cells.Item(3028, 5).Value = 221
$ws.Cells.Item(3029, 1).Value = 44174
$ws.Cells.Item(3029, 2).Value = "Okres Svidník"
$ws.Cells.Item(3029, 3).Value = 145
$ws.Cells.Item(3029, 4).Value = 10
$ws.Cells.Item(3029, 5).Value = 155
$ws.Cells.Item(3030, 1).Value = 44174
$ws.Cells.Item(3030, 2).Value = "Okres Topoľčany"
$ws.Cells.Item(3030, 3).Value = 996
$ws.Cells.Item(3030, 4).Value = 90
$ws.Cells.Item(3030, 5).Value = 1086
$ws.Cells.Item(3031, 1).Value = 44174
$ws.Cells.Item(3031, 2).Value = "Okres Trebišov"
$ws.Cells.Item(3031, 3).Value = 304
$ws.Cells.Item(3031, 4).Value = 13
$ws.Cells.Item(3031, 5).Value = 317
$ws.Cells.Item(3032, 1).Value = 44174
$ws.Cells.Item(3032, 2).Value = "Okres Trenčín"
$ws.Cells.Item(3032, 3).Value = 499
$ws.Cells.Item(3032, 4).Value = 59
$ws.Cells.Item(3032, 5).Value = 558
$ws.Cells.Item(3033, 1).Value = 44174
$ws.Cells.Item(3033, 2).Value = "Okres Trnava"
$ws.Cells.Item(3033, 3).Value = 510
$ws.Cells.Item(3033, 4).Value = 47
$ws.Cells.Item(3033, 5).Value = 557
$ws.Cells.Item(3034, 1).Value = 44174
$ws.Cells.Item(3034, 2).Value = "Okres Tvrdošín"
$ws.Cells.Item(3034, 3).Value = 194
$ws.Cells.Item(3034, 4).Value = 17
$ws.Cells.Item(3034, 5).Value = 211
$ws.Cells.Item(3035, 1).Value = 44174
$ws.Cells.Item(3035, 2).Value = "Okres Vranov nad Topľou"
$ws.Cells.Item(3035, 3).Value = 274
$ws.Cells.Item(3035, 4).Value = 24
$ws.Cells.Item(3035, 5).Value = 298
$ws.Cells.Item(3036, 1).Value = 44174
$ws.Cells.Item(3036, 2).Value = "Okres Zlaté Moravce"
$ws.Cells.Item(3036, 3).Value = 2
$ws.Cells.Item(3036, 4).Value = 0
$ws.Cells.Item(3036, 5).Value = 2
$ws.Cells.Item(3037, 1).Value = 44174
$ws.Cells.Item(3037, 2).Value = "Okres Zvolen"
$ws.Cells.Item(3037, 3).Value = 213
$ws.Cells.Item(3037, 4).Value = 15
$ws.Cells.Item(3037, 5).Value = 228
$ws.Cells.Item(3038, 1).Value = 44174
$ws.Cells.Item(3038, 2).Value = "Okres Žarnovica"
$ws.Cells.Item(3038, 3).Value = 792
$ws.Cells.Item(3038, 4).Value = 69
$ws.Cells.Item(3038, 5).Value = 861
$ws.Cells.Item(3039, 1).Value = 44174
$ws.Cells.Item(3039, 2).Value = "Okres Žiar nad Hronom"
$ws.Cells.Item(3039, 3).Value = 356
$ws.Cells.Item(3039, 4).Value = 22
$ws.Cells.Item(3039, 5).Value = 378
$ws.Cells.Item(3040, 1).Value = 44174
$ws.Cells.Item(3040, 2).Value = "Okres Žilina"
$ws.Cells.Item(3040, 3).Value = 179
$ws.Cells.Item(3040, 4).Value = 13
$ws.Cells.Item(3040, 5).Value = 192
$ws.Cells.Item(3041, 1).Value = 44175
$ws.Cells.Item(3041, 2).Value = "Okres Banská Bystrica"
$ws.Cells.Item(3041, 3).Value = 1122
$ws.Cells.Item(3041, 4).Value = 51
$ws.Cells.Item(3041, 5).Value = 1173
$ws.Cells.Item(3042, 1).Value = 44175
$ws.Cells.Item(3042, 2).Value = "Okres Bardejov"
$ws.Cells.Item(3042, 3).Value = 1
$ws.Cells.Item(3042, 4).Value = 0
$ws.Cells.Item(3042, 5).Value = 1
$ws.Cells.Item(3043, 1).Value = 44175
$ws.Cells.Item(3043, 2).Value = "Okres Bratislava I"
$ws.Cells.Item(3043, 3).Value = 75
$ws.Cells.Item(3043, 4).Value = 0
$ws.Cells.Item(3043, 5).Value = 75
$ws.Cells.Item(3044, 1).Value = 44175
$ws.Cells.Item(3044, 2).Value = "Okres Bratislava II"
$ws.Cells.Item(3044, 3).Value = 680
$ws.Cells.Item(3044, 4).Value = 34
$ws.Cells.Item(3044, 5).Value = 714
$ws.Cells.Item(3045, 1).Value = 44175
$ws.Cells.Item(3045, 2).Value = "Okres Bratislava III"
$ws.Cells.Item(3045, 3).Value = 683
$ws.Cells.Item(3045, 4).Value = 24
$ws.Cells.Item(3045, 5).Value = 707
$ws.Cells.Item(3046, 1).Value = 44175
$ws.Cells.Item(3046, 2).Value = "Okres Bratislava V"
$ws.Cells.Item(3046, 3).Value = 266
$ws.Cells.Item(3046, 4).Value = 12
$ws.Cells.Item(3046, 5).Value = 278
$ws.Cells.Item(3047, 1).Value = 44175
$ws.Cells.Item(3047, 2).Value = "Okres Brezno"
$ws.Cells.Item(3047, 3).Value = 429
$ws.Cells.Item(3047, 4).Value = 116
$ws.Cells.Item(3047, 5).Value = 545
$ws.Cells.Item(3048, 1).Value = 44175
$ws.Cells.Item(3048, 2).Value = "Okres Čadca"
$ws.Cells.Item(3048, 3).Value = 432
$ws.Cells.Item(3048, 4).Value = 29
$ws.Cells.Item(3048, 5).Value = 461
$ws.Cells.Item(3049, 1).Value = 44175
$ws.Cells.Item(3049, 2).Value = "Okres Dolný Kubín"
$ws.Cells.Item(3049, 3).Value = 35
$ws.Cells.Item(3049, 4).Value = 144
$ws.Cells.Item(3049, 5).Value = 179
$ws.Cells.Item(3050, 1).Value = 44175
$ws.Cells.Item(3050, 2).Value = "Okres Dunajská Streda"
$ws.Cells.Item(3050, 3).Value = 355
$ws.Cells.Item(3050, 4).Value = 18
$ws.Cells.Item(3050, 5).Value = 373
$ws.Cells.Item(3051, 1).Value = 44175
$ws.Cells.Item(3051, 2).Value = "Okres Galanta"
$ws.Cells.Item(3051, 3).Value = 574
$ws.Cells.Item(3051, 4).Value = 32
$ws.Cells.Item(3051, 5).Value = 606
$ws.Cells.Item(3052, 1).Value = 44175
$ws.Cells.Item(3052, 2).Value = "Okres Gelnica"
$ws.Cells.Item(3052, 3).Value = 214
$ws.Cells.Item(3052, 4).Value = 8
$ws.Cells.Item(3052, 5).Value = 222
$ws.Cells.Item(3053, 1).Value = 44175
$ws.Cells.Item(3053, 2).Value = "Okres Hlohovec"
$ws.Cells.Item(3053, 3).Value = 19
$ws.Cells.Item(3053, 4).Value = 2
$ws.Cells.Item(3053, 5).Value = 21
$ws.Cells.Item(3054, 1).Value = 44175
$ws.Cells.Item(3054, 2).Value = "Okres Humenné"
$ws.Cells.Item(3054, 3).Value = 231
$ws.Cells.Item(3054, 4).Value = 10
$ws.Cells.Item(3054, 5).Value = 241
$ws.Cells.Item(3055, 1).Value = 44175
$ws.Cells.Item(3055, 2).Value = "Okres Ilava"
$ws.Cells.Item(3055, 3).Value = 355
$ws.Cells.Item(3055, 4).Value = 44
$ws.Cells.Item(3055, 5).Value = 399
$ws.Cells.Item(3056, 1).Value = 44175
$ws.Cells.Item(3056, 2).Value = "Okres Kežmarok"
$ws.Cells.Item(3056, 3).Value = 220
$ws.Cells.Item(3056, 4).Value = 16
$ws.Cells.Item(3056, 5).Value = 236
$ws.Cells.Item(3057, 1).Value = 44175
$ws.Cells.Item(3057, 2).Value = "Okres Komárno"
$ws.Cells.Item(3057, 3).Value = 363
$ws.Cells.Item(3057, 4).Value = 12
$ws.Cells.Item(3057, 5).Value = 375
$ws.Cells.Item(3058, 1).Value = 44175
$ws.Cells.Item(3058, 2).Value = "Okres Košice I"
$ws.Cells.Item(3058, 3).Value = 453
$ws.Cells.Item(3058, 4).Value = 33
$ws.Cells.Item(3058, 5).Value = 486
$ws.Cells.Item(3059, 1).Value = 44175
$ws.Cells.Item(3059, 2).Value = "Okres Košice II"
$ws.Cells.Item(3059, 3).Value = 619
$ws.Cells.Item(3059, 4).Value = 15
$ws.Cells.Item(3059, 5).Value = 634
$ws.Cells.Item(3060, 1).Value = 44175
$ws.Cells.Item(3060, 2).Value = "Okres Košice IV"
$ws.Cells.Item(3060, 3).Value = 198
$ws.Cells.Item(3060, 4).Value = 7
$ws.Cells.Item(3060, 5).Value = 205
$ws.Cells.Item(3061, 1).Value = 44175
$ws.Cells.Item(3061, 2).Value = "Okres Krupina"
$ws.Cells.Item(3061, 3).Value = 2
$ws.Cells.Item(3061, 4).Value = 1
$ws.Cells.Item(3061, 5).Value = 3
$ws.Cells.Item(3062, 1).Value = 44175
$ws.Cells.Item(3062, 2).Value = "Okres Levice"
$ws.Cells.Item(3062, 3).Value = 500
$ws.Cells.Item(3062, 4).Value = 34
$ws.Cells.Item(3062, 5).Value = 534
$ws.Cells.Item(3063, 1).Value = 44175
$ws.Cells.Item(3063, 2).Value = "Okres Levoča"
$ws.Cells.Item(3063, 3).Value = 138
$ws.Cells.Item(3063, 4).Value = 13
$ws.Cells.Item(3063, 5).Value = 151
$ws.Cells.Item(3064, 1).Value = 44175
$ws.Cells.Item(3064, 2).Value = "Okres Liptovský Mikuláš"
$ws.Cells.Item(3064, 3).Value = 217
$ws.Cells.Item(3064, 4).Value = 11
$ws.Cells.Item(3064, 5).Value = 228
$ws.Cells.Item(3065, 1).Value = 44175
$ws.Cells.Item(3065, 2).Value = "Okres Lučenec"
$ws.Cells.Item(3065, 3).Value = 295
$ws.Cells.Item(3065, 4).Value = 16
$ws.Cells.Item(3065, 5).Value = 311
$ws.Cells.Item(3066, 1).Value = 44175
$ws.Cells.Item(3066, 2).Value = "Okres Malacky"
$ws.Cells.Item(3066, 3).Value = 190
$ws.Cells.Item(3066, 4).Value = 15
$ws.Cells.Item(3066, 5).Value = 205
$ws.Cells.Item(3067, 1).Value = 44175
$ws.Cells.Item(3067, 2).Value = "Okres Martin"
$ws.Cells.Item(3067, 3).Value = 506
$ws.Cells.Item(3067, 4).Value = 42
$ws.Cells.Item(3067, 5).Value = 548
$ws.Cells.Item(3068, 1).Value = 44175
$ws.Cells.Item(3068, 2).Value = "Okres Michalovce"
$ws.Cells.Item(3068, 3).Value = 552
$ws.Cells.Item(3068, 4).Value = 42
$ws.Cells.Item(3068, 5).Value = 594
$ws.Cells.Item(3069, 1).Value = 44175
$ws.Cells.Item(3069, 2).Value = "Okres Myjava"
$ws.Cells.Item(3069, 3).Value = 195
$ws.Cells.Item(3069, 4).Value = 29
$ws.Cells.Item(3069, 5).Value = 224
$ws.Cells.Item(3070, 1).Value = 44175
$ws.Cells.Item(3070, 2).Value = "Okres Nitra"
$ws.Cells.Item(3070, 3).Value = 602
$ws.Cells.Item(3070, 4).Value = 44
$ws.Cells.Item(3070, 5).Value = 646
$ws.Cells.Item(3071, 1).Value = 44175
$ws.Cells.Item(3071, 2).Value = "Okres Nové Mesto nad Váhom"
$ws.Cells.Item(3071, 3).Value = 188
$ws.Cells.Item(3071, 4).Value = 27
$ws.Cells.Item(3071, 5).Value = 215
$ws.Cells.Item(3072, 1).Value = 44175
$ws.Cells.Item(3072, 2).Value = "Okres Partizánske"
$ws.Cells.Item(3072, 3).Value = 204
$ws.Cells.Item(3072, 4).Value = 19
$ws.Cells.Item(3072, 5).Value = 223
$ws.Cells.Item(3073, 1).Value = 44175
$ws.Cells.Item(3073, 2).Value = "Okres Pezinok"
$ws.Cells.Item(3073, 3).Value = 376
$ws.Cells.Item(3073, 4).Value = 12
$ws.Cells.Item(3073, 5).Value = 388
$ws.Cells.Item(3074, 1).Value = 44175
$ws.Cells.Item(3074, 2).Value = "Okres Piešťany"
$ws.Cells.Item(3074, 3).Value = 144
$ws.Cells.Item(3074, 4).Value = 5
$ws.Cells.Item(3074, 5).Value = 149
$ws.Cells.Item(3075, 1).Value = 44175
$ws.Cells.Item(3075, 2).Value = "Okres Poprad"
$ws.Cells.Item(3075, 3).Value = 948
$ws.Cells.Item(3075, 4).Value = 113
$ws.Cells.Item(3075, 5).Value = 1061
$ws.Cells.Item(3076, 1).Value = 44175
$ws.Cells.Item(3076, 2).Value = "Okres Považská Bystrica"
$ws.Cells.Item(3076, 3).Value = 290
$ws.Cells.Item(3076, 4).Value = 39
$ws.Cells.Item(3076, 5).Value = 329
$ws.Cells.Item(3077, 1).Value = 44175
$ws.Cells.Item(3077, 2).Value = "Okres Prešov"
$ws.Cells.Item(3077, 3).Value = 233
$ws.Cells.Item(3077, 4).Value = 11
$ws.Cells.Item(3077, 5).Value = 244
$ws.Cells.Item(3078, 1).Value = 44175
$ws.Cells.Item(3078, 2).Value = "Okres Prievidza"
$ws.Cells.Item(3078, 3).Value = 393
$ws.Cells.Item(3078, 4).Value = 39
$ws.Cells.Item(3078, 5).Value = 432
$ws.Cells.Item(3079, 1).Value = 44175
$ws.Cells.Item(3079, 2).Value = "Okres Púchov"
$ws.Cells.Item(3079, 3).Value = 4
$ws.Cells.Item(3079, 4).Value = 1
$ws.Cells.Item(3079, 5).Value = 5
$ws.Cells.Item(3080, 1).Value = 44175
$ws.Cells.Item(3080, 2).Value = "Okres Revúca"
$ws.Cells.Item(3080, 3).Value = 259
$ws.Cells.Item(3080, 4).Value = 5
$ws.Cells.Item(3080, 5).Value = 264
$ws.Cells.Item(3081, 1).Value = 44175
$ws.Cells.Item(3081, 2).Value = "Okres Rimavská Sobota"
$ws.Cells.Item(3081, 3).Value = 241
$ws.Cells.Item(3081, 4).Value = 10
$ws.Cells.Item(3081, 5).Value = 251
$ws.Cells.Item(3082, 1).Value = 44175
$ws.Cells.Item(3082, 2).Value = "Okres Rožňava"
$ws.Cells.Item(3082, 3).Value = 224
$ws.Cells.Item(3082, 4).Value = 3
$ws.Cells.Item(3082, 5).Value = 227
$ws.Cells.Item(3083, 1).Value = 44175
$ws.Cells.Item(3083, 2).Value = "Okres Ružomberok"
$ws.Cells.Item(3083, 3).Value = 240
$ws.Cells.Item(3083, 4).Value = 17
$ws.Cells.Item(3083, 5).Value = 257
$ws.Cells.Item(3084, 1).Value = 44175
$ws.Cells.Item(3084, 2).Value = "Okres Skalica"
$ws.Cells.Item(3084, 3).Value = 563
$ws.Cells.Item(3084, 4).Value = 16
$ws.Cells.Item(3084, 5).Value = 579
$ws.Cells.Item(3085, 1).Value = 44175
$ws.Cells.Item(3085, 2).Value = "Okres Snina"
$ws.Cells.Item(3085, 3).Value = 108
$ws.Cells.Item(3085, 4).Value = 11
$ws.Cells.Item(3085, 5).Value = 119
$ws.Cells.Item(3086, 1).Value = 44175
$ws.Cells.Item(3086, 2).Value = "Okres Sobrance"
$ws.Cells.Item(3086, 3).Value = 107
$ws.Cells.Item(3086, 4).Value = 4
$ws.Cells.Item(3086, 5).Value = 111
$ws.Cells.Item(3087, 1).Value = 44175
$ws.Cells.Item(3087, 2).Value = "Okres Spišská Nová Ves"
$ws.Cells.Item(3087, 3).Value = 130
$ws.Cells.Item(3087, 4).Value = 14
$ws.Cells.Item(3087, 5).Value = 144
$ws.Cells.Item(3088, 1).Value = 44175
$ws.Cells.Item(3088, 2).Value = "Okres Stará Ľubovňa"
$ws.Cells.Item(3088, 3).Value = 180
$ws.Cells.Item(3088, 4).Value = 21
$ws.Cells.Item(3088, 5).Value = 201
$ws.Cells.Item(3089, 1).Value = 44175
$ws.Cells.Item(3089, 2).Value = "Okres Svidník"
$ws.Cells.Item(3089, 3).Value = 39
$ws.Cells.Item(3089, 4).Value = 2
$ws.Cells.Item(3089, 5).Value = 41
$ws.Cells.Item(3090, 1).Value = 44175
$ws.Cells.Item(3090, 2).Value = "Okres Topoľčany"
$ws.Cells.Item(3090, 3).Value = 181
$ws.Cells.Item(3090, 4).Value = 9
$ws.Cells.Item(3090, 5).Value = 190
$ws.Cells.Item(3091, 1).Value = 44175
$ws.Cells.Item(3091, 2).Value = "Okres Trebišov"
$ws.Cells.Item(3091, 3).Value = 628
$ws.Cells.Item(3091, 4).Value = 11
$ws.Cells.Item(3091, 5).Value = 639
$ws.Cells.Item(3092, 1).Value = 44175
$ws.Cells.Item(3092, 2).Value = "Okres Trenčín"
$ws.Cells.Item(3092, 3).Value = 275
$ws.Cells.Item(3092, 4).Value = 37
$ws.Cells.Item(3092, 5).Value = 312
$ws.Cells.Item(3093, 1).Value = 44175
$ws.Cells.Item(3093, 2).Value = "Okres Trnava"
$ws.Cells.Item(3093, 3).Value = 603
$ws.Cells.Item(3093, 4).Value = 30
$ws.Cells.Item(3093, 5).Value = 633
$ws.Cells.Item(3094, 1).Value = 44175
$ws.Cells.Item(3094, 2).Value = "Okres Tvrdošín"
$ws.Cells.Item(3094, 3).Value = 45
$ws.Cells.Item(3094, 4).Value = 6
$ws.Cells.Item(3094, 5).Value = 51
$ws.Cells.Item(3095, 1).Value = 44175
$ws.Cells.Item(3095, 2).Value = "Okres Veľký Krtíš"
$ws.Cells.Item(3095, 3).Value = 342
$ws.Cells.Item(3095, 4).Value = 9
$ws.Cells.Item(3095, 5).Value = 351
$ws.Cells.Item(3096, 1).Value = 44175
$ws.Cells.Item(3096, 2).Value = "Okres Vranov nad Topľou"
$ws.Cells.Item(3096, 3).Value = 338
$ws.Cells.Item(3096, 4).Value = 32
$ws.Cells.Item(3096, 5).Value = 370
$ws.Cells.Item(3097, 1).Value = 44175
$ws.Cells.Item(3097, 2).Value = "Okres Zlaté Moravce"
$ws.Cells.Item(3097, 3).Value = 2
$ws.Cells.Item(3097, 4).Value = 0
$ws.Cells.Item(3097, 5).Value = 2
$ws.Cells.Item(3098, 1).Value = 44175
$ws.Cells.Item(3098, 2).Value = "Okres Zvolen"
$ws.Cells.Item(3098, 3).Value = 186
$ws.Cells.Item(3098, 4).Value = 17
$ws.Cells.Item(3098, 5).Value = 203
$ws.Cells.Item(3099, 1).Value = 44175
$ws.Cells.Item(3099, 2).Value = "Okres Žiar nad Hronom"
$ws.Cells.Item(3099, 3).Value = 710
$ws.Cells.Item(3099, 4).Value = 31
$ws.Cells.Item(3099, 5).Value = 741
$ws.Cells.Item(3100, 1).Value = 44175
$ws.Cells.Item(3100, 2).Value = "Okres Žilina"
$ws.Cells.Item(3100, 3).Value = 6740
$ws.Cells.Item(3100, 4).Value = 689
$ws.Cells.Item(3100, 5).Value = 7429
